$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 83, shifting existing rows 83-92 down to 84-93.
$ws.Rows("83:83").Insert()

# Copy the date number-format style used by column D (style "s=2" -> numFmt 165)
# from the row above so the new row's date cell is formatted the same way.
$ws.Range("D82").Copy()
$ws.Range("D83").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 83 with the weekly record values.
$ws.Range("A83").Value = 5
$ws.Range("B83").Value = "Macroferia Regional de Talca"
$ws.Range("C83").Value = "Maule"
$ws.Range("D83").Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E83").Value = 7
$ws.Range("F83").Value = 100112022
$ws.Range("G83").Value = "Arveja Verde"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 200
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 27000
$ws.Range("M83").Value = 27000
$ws.Range("N83").Value = "`$/saco 25 kilos"
$ws.Range("O83").Value = "Carahue"
$ws.Range("P83").Value = 1080
$ws.Range("Q83").Value = 25
$ws.Range("R83").Value = "Hortaliza"
